# Generate Report for Handback
# Applies the "handback" status update to the localization-status workbook:
#  - Overview sheet: status cells flip from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: populate the "Latest Target File" / "Latest
#    Handback File" / "Latest Handback DateTime" columns (I/J/K) for both
#    data rows, including hyperlinks on column I mirroring column A's link,
#    and widen a few columns that now hold longer filenames.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ColumnWidth inputs chosen so the engine's internal (chars -> stored width)
# rounding lands on the desired stored widths (30 and 40 character units).
$wWide = 29.16796875   # -> stored width 30 (widens the "~17.2" columns)
$wFull = 39.16796875   # -> stored width 40 (widens the "~18.6"/"~21.7" columns)

# ---- Column widths -------------------------------------------------------

$overview.Range("E1").ColumnWidth = $wWide
$overview.Range("F1").ColumnWidth = $wWide

$zhcn.Range("C1").ColumnWidth = $wWide
$zhcn.Range("I1").ColumnWidth = $wFull
$zhcn.Range("J1").ColumnWidth = $wFull

$dede.Range("C1").ColumnWidth = $wWide
$dede.Range("I1").ColumnWidth = $wFull
$dede.Range("J1").ColumnWidth = $wFull

# ---- Status text ----------------------------------------------------------
# All eight cells below shared the single "Ready for handoff" string in the
# source workbook (Overview's per-language status columns AND each language
# sheet's own "Status" column for both rows) - the handback flips all of
# them to the same new status text.

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# ---- Hyperlink targets (same URLs column A already links to) ------------

$target1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e69c6f118dcc65e68dc1cb49b3f8b20279af71e9/e2e/1e15f338-1cba-4085-984c-8b2cf1d82c21.md"
$target2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e69c6f118dcc65e68dc1cb49b3f8b20279af71e9/e2e/a3f895e8-ad94-4c6a-ab44-2e29423c8cae.md"

# ---- zh-cn: Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K)

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $target1, "", "", "1e15f338-1cba-4085-984c-8b2cf1d82c21.md")
$zhcn.Range("J2").Value = "1e15f338-1cba-4085-984c-8b2cf1d82c21.2ca43822db748e0249d3fcc718f6e3c6620efced.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-24 11:04:28"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $target2, "", "", "a3f895e8-ad94-4c6a-ab44-2e29423c8cae.md")
$zhcn.Range("J3").Value = "a3f895e8-ad94-4c6a-ab44-2e29423c8cae.b46cd5dba1cbd5abdc5de711928ef2a9fc6dce02.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-24 11:04:28"

# ---- de-de: Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K)

$dede.Hyperlinks.Add($dede.Range("I2"), $target1, "", "", "1e15f338-1cba-4085-984c-8b2cf1d82c21.md")
$dede.Range("J2").Value = "1e15f338-1cba-4085-984c-8b2cf1d82c21.2ca43822db748e0249d3fcc718f6e3c6620efced.de-de.xlf"
$dede.Range("K2").Value = "2016-08-24 11:04:35"

$dede.Hyperlinks.Add($dede.Range("I3"), $target2, "", "", "a3f895e8-ad94-4c6a-ab44-2e29423c8cae.md")
$dede.Range("J3").Value = "a3f895e8-ad94-4c6a-ab44-2e29423c8cae.b46cd5dba1cbd5abdc5de711928ef2a9fc6dce02.de-de.xlf"
$dede.Range("K3").Value = "2016-08-24 11:04:35"
